$wb = $excel.ActiveWorkbook

# ===== Sheet ALC =====
$ws = $wb.Worksheets.Item("ALC")
# hunk originally near old-file line 4342
$ws.Range("H74").Value = 5900.375
$ws.Range("I74").Value = 4900.6
$ws.Range("J74").Value = 7566.6665
$ws.Range("K74").Value = 4900.6
$ws.Range("L74").Value = 7566.6665
$ws.Range("M74").Value = -3964.6
$ws.Range("N74").Value = -9438.666499999999

# hunk originally near old-file line 4492
$ws.Range("H77").Value = 5900.375
$ws.Range("I77").Value = 4900.6
$ws.Range("J77").Value = 7566.6665
$ws.Range("K77").Value = 24503
$ws.Range("L77").Value = 37833.3325
$ws.Range("M77").Value = -19823
$ws.Range("N77").Value = -47193.3325

# hunk originally near old-file line 5811
$ws.Range("H103").Value = 691.5161000000001
$ws.Range("I103").Value = 400
$ws.Range("J103").Value = 734.7037
$ws.Range("K103").Value = 1200
$ws.Range("L103").Value = 2204.1111
$ws.Range("M103").Value = -614
$ws.Range("N103").Value = -3376.1111

# hunk originally near old-file line 6408
$ws.Range("H115").Value = 655.7143
$ws.Range("I115").Value = 598.3333
$ws.Range("J115").Value = 1000
$ws.Range("K115").Value = 1794.9999
$ws.Range("L115").Value = 3000
$ws.Range("M115").Value = -227.9999
$ws.Range("N115").Value = -6134

# hunk originally near old-file line 7504
$ws.Range("H137").Value = 1630.4572
$ws.Range("I137").Value = 1266.6666
$ws.Range("J137").Value = 1664.5625
$ws.Range("K137").Value = 3799.9998
$ws.Range("L137").Value = 4993.6875
$ws.Range("M137").Value = -1249.9998
$ws.Range("N137").Value = -10093.6875


# ===== Sheet ARM =====
$ws = $wb.Worksheets.Item("ARM")
# hunk originally near old-file line 9974
$ws.Range("H45").Value = 1433.2778
$ws.Range("I45").Value = 1347.1818
$ws.Range("J45").Value = 1568.5714
$ws.Range("K45").Value = 1347.1818
$ws.Range("L45").Value = 1568.5714
$ws.Range("M45").Value = -970.1818000000001
$ws.Range("N45").Value = -2322.5714

# hunk originally near old-file line 10740
$ws.Range("H61").Value = 1546.6364
$ws.Range("I61").Value = 969.4194
$ws.Range("J61").Value = 2923.077
$ws.Range("K61").Value = 969.4194
$ws.Range("L61").Value = 2923.077
$ws.Range("M61").Value = -757.4194
$ws.Range("N61").Value = -3347.077

# hunk originally near old-file line 11374
$ws.Range("H74").Value = 1592.2333
$ws.Range("I74").Value = 843.4167
$ws.Range("J74").Value = 2091.4443
$ws.Range("K74").Value = 843.4167
$ws.Range("L74").Value = 2091.4443
$ws.Range("M74").Value = 30.58330000000001
$ws.Range("N74").Value = -3839.4443

# hunk originally near old-file line 11524
$ws.Range("H77").Value = 1592.2333
$ws.Range("I77").Value = 843.4167
$ws.Range("J77").Value = 2091.4443
$ws.Range("K77").Value = 4217.0835
$ws.Range("L77").Value = 10457.2215
$ws.Range("M77").Value = 150.9165000000003
$ws.Range("N77").Value = -19193.2215

# hunk originally near old-file line 12743
$ws.Range("H102").Value = 93241.73
$ws.Range("I102").Value = 144997
$ws.Range("K102").Value = 144997
$ws.Range("M102").Value = -143375

# hunk originally near old-file line 14210
$ws.Range("H132").Value = 2119.3157
$ws.Range("I132").Value = 1963.0968
$ws.Range("J132").Value = 2811.1428
$ws.Range("K132").Value = 5889.2904
$ws.Range("L132").Value = 8433.428400000001
$ws.Range("M132").Value = -3359.2904
$ws.Range("N132").Value = -13493.4284

# hunk originally near old-file line 14409
$ws.Range("H136").Value = 1546.6364
$ws.Range("I136").Value = 969.4194
$ws.Range("J136").Value = 2923.077
$ws.Range("K136").Value = 2908.2582
$ws.Range("L136").Value = 8769.231
$ws.Range("M136").Value = -358.2582000000002
$ws.Range("N136").Value = -13869.231


# ===== Sheet BSM =====
$ws = $wb.Worksheets.Item("BSM")
# hunk originally near old-file line 15051
$ws.Range("H7").Value = 250475
$ws.Range("I7").Value = 250475
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 250475
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -250362
$ws.Range("N7").ClearContents()

# hunk originally near old-file line 19844
$ws.Range("H105").Value = 183981.81
$ws.Range("I105").Value = 144655.58
$ws.Range("K105").Value = 144655.58
$ws.Range("M105").Value = -142908.58


# ===== Sheet CRP =====
$ws = $wb.Worksheets.Item("CRP")
# hunk originally near old-file line 21739
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

# hunk originally near old-file line 23166
$ws.Range("H31").Value = 1304.725
$ws.Range("I31").Value = 945.53845
$ws.Range("J31").Value = 1646.3903
$ws.Range("K31").Value = 945.53845
$ws.Range("L31").Value = 1646.3903
$ws.Range("M31").Value = -650.53845
$ws.Range("N31").Value = -2236.3903

# hunk originally near old-file line 23316
$ws.Range("H34").Value = 1304.725
$ws.Range("I34").Value = 945.53845
$ws.Range("J34").Value = 1646.3903
$ws.Range("K34").Value = 945.53845
$ws.Range("L34").Value = 1646.3903
$ws.Range("M34").Value = -743.53845
$ws.Range("N34").Value = -2050.3903

# hunk originally near old-file line 24486
$ws.Range("H58").Value = 1310.4222
$ws.Range("I58").Value = 1152.9286
$ws.Range("J58").Value = 1569.8235
$ws.Range("K58").Value = 1152.9286
$ws.Range("L58").Value = 1569.8235
$ws.Range("M58").Value = -949.9286
$ws.Range("N58").Value = -1975.8235

# hunk originally near old-file line 24691
$ws.Range("H62").Value = 4110.4443
$ws.Range("I62").Value = 2628.5
$ws.Range("J62").Value = 5296
$ws.Range("K62").Value = 2628.5
$ws.Range("L62").Value = 5296
$ws.Range("M62").Value = -2004.5
$ws.Range("N62").Value = -6544

# hunk originally near old-file line 24841
$ws.Range("H65").Value = 4110.4443
$ws.Range("I65").Value = 2628.5
$ws.Range("J65").Value = 5296
$ws.Range("K65").Value = 13142.5
$ws.Range("L65").Value = 26480
$ws.Range("M65").Value = -10022.5
$ws.Range("N65").Value = -32720

# hunk originally near old-file line 26489
$ws.Range("H99").Value = 14624.375
$ws.Range("I99").Value = 2264
$ws.Range("J99").Value = 22040.6
$ws.Range("K99").Value = 2264
$ws.Range("L99").Value = 22040.6
$ws.Range("M99").Value = -766
$ws.Range("N99").Value = -25036.6

# hunk originally near old-file line 27619
$ws.Range("H122").Value = 1965.0714
$ws.Range("I122").Value = 2377.5
$ws.Range("J122").Value = 1415.1666
$ws.Range("K122").Value = 7132.5
$ws.Range("L122").Value = 4245.4998
$ws.Range("M122").Value = -4682.5
$ws.Range("N122").Value = -9145.4998

# hunk originally near old-file line 27818
$ws.Range("H126").Value = 14624.375
$ws.Range("I126").Value = 2264
$ws.Range("J126").Value = 22040.6
$ws.Range("K126").Value = 6792
$ws.Range("L126").Value = 66121.79999999999
$ws.Range("M126").Value = -4322
$ws.Range("N126").Value = -71061.79999999999

# hunk originally near old-file line 28109
$ws.Range("H132").Value = 3031
$ws.Range("I132").Value = 2896.5833
$ws.Range("K132").Value = 8689.749899999999
$ws.Range("M132").Value = -6159.749899999999

# hunk originally near old-file line 28210
$ws.Range("H134").Value = 1655.8462
$ws.Range("I134").Value = 1228.4546
$ws.Range("K134").Value = 3685.3638
$ws.Range("M134").Value = -1150.3638

# hunk originally near old-file line 28311
$ws.Range("H136").Value = 1310.4222
$ws.Range("I136").Value = 1152.9286
$ws.Range("J136").Value = 1569.8235
$ws.Range("K136").Value = 3458.7858
$ws.Range("L136").Value = 4709.470499999999
$ws.Range("M136").Value = -908.7857999999997
$ws.Range("N136").Value = -9809.470499999999


# ===== Sheet GSM =====
$ws = $wb.Worksheets.Item("GSM")
# hunk originally near old-file line 41281
$ws.Range("H113").Value = 1423.375
$ws.Range("I113").Value = 1166.5834
$ws.Range("K113").Value = 1166.5834
$ws.Range("M113").Value = 1003.4166

# hunk originally near old-file line 41915
$ws.Range("H126").Value = 3003.0908
$ws.Range("I126").Value = 3294.25
$ws.Range("J126").Value = 2226.6667
$ws.Range("K126").Value = 9882.75
$ws.Range("L126").Value = 6680.000100000001
$ws.Range("M126").Value = -7412.75
$ws.Range("N126").Value = -11620.0001

# hunk originally near old-file line 42212
$ws.Range("H132").Value = 2126.5667
$ws.Range("I132").Value = 1874.037
$ws.Range("J132").Value = 4399.3335
$ws.Range("K132").Value = 5622.111
$ws.Range("L132").Value = 13198.0005
$ws.Range("M132").Value = -3092.111
$ws.Range("N132").Value = -18258.0005


# ===== Sheet LTW =====
$ws = $wb.Worksheets.Item("LTW")
# hunk originally near old-file line 47243
$ws.Range("H93").Value = 1365.6154
$ws.Range("I93").Value = 1341.1818
$ws.Range("K93").Value = 1341.1818
$ws.Range("M93").Value = -93.18180000000007

# hunk originally near old-file line 48661
$ws.Range("H122").Value = 2545.818
$ws.Range("I122").Value = 2600.4
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 7801.200000000001
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -5351.200000000001
$ws.Range("N122").Value = -10900

# hunk originally near old-file line 49151
$ws.Range("H132").Value = 5257.5713
$ws.Range("I132").Value = 5805.778
$ws.Range("K132").Value = 17417.334
$ws.Range("M132").Value = -14887.334

# hunk originally near old-file line 49350
$ws.Range("H136").Value = 1078.7675
$ws.Range("I136").Value = 966.7241
$ws.Range("J136").Value = 1310.8572
$ws.Range("K136").Value = 2900.1723
$ws.Range("L136").Value = 3932.5716
$ws.Range("M136").Value = -350.1723000000002
$ws.Range("N136").Value = -9032.571599999999


# ===== Sheet WVR =====
$ws = $wb.Worksheets.Item("WVR")
# hunk originally near old-file line 54386
$ws.Range("H96").Value = 125003140
$ws.Range("I96").Value = 200003170
$ws.Range("J96").Value = 3093.3333
$ws.Range("K96").Value = 200003170
$ws.Range("L96").Value = 3093.3333
$ws.Range("M96").Value = -200001797
$ws.Range("N96").Value = -5839.3333

# hunk originally near old-file line 55859
$ws.Range("H126").Value = 1924.25
$ws.Range("I126").Value = 1758.3
$ws.Range("J126").Value = 2200.8333
$ws.Range("K126").Value = 5274.9
$ws.Range("L126").Value = 6602.499899999999
$ws.Range("M126").Value = -2804.9
$ws.Range("N126").Value = -11542.4999

